$p = $ppt.ActivePresentation

# Locate the "TextBox 77" shape on slide 1 that holds the
# "deletePerson(p)" sequence-diagram call-out, and rename the
# method call to "deleteTask(t)" while preserving the existing
# run/formatting split (the label is made of two runs:
# "deletePerson" + "(p)").
$s = $ppt.ActivePresentation.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "deletePerson(p)") {
                $target = $sh
                break
            }
        }
    }
}

$tr = $target.TextFrame2.TextRange

# First run: "deletePerson" -> "deleteTask"
$tr.Characters(1, 12).Text = "deleteTask"

# Second run: "(p)" -> "(t)" (now starts right after "deleteTask")
$tr.Characters(11, 3).Text = "(t)"
